$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 4651.5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 4651.5
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 13954.5
$ws.Range("N17").Value = -14290.5

$ws.Range("H32").Value = 10965.223
$ws.Range("I32").Value = 9964.666999999999
$ws.Range("J32").Value = 11465.5
$ws.Range("K32").Value = 9964.666999999999
$ws.Range("L32").Value = 11465.5
$ws.Range("M32").Value = -9638.666999999999

$ws.Range("H42").Value = 1185.75
$ws.Range("I42").Value = 93.75
$ws.Range("J42").Value = 3369.75
$ws.Range("K42").Value = 281.25
$ws.Range("L42").Value = 10109.25
$ws.Range("M42").Value = -51.25
$ws.Range("N42").Value = -10569.25

$ws.Range("H51").Value = 12276.211
$ws.Range("I51").Value = 13665.333
$ws.Range("J51").Value = 11026
$ws.Range("K51").Value = 13665.333
$ws.Range("L51").Value = 11026
$ws.Range("M51").Value = -13181.333
$ws.Range("N51").Value = -11994

$ws.Range("H112").Value = 6012.9165
$ws.Range("I112").Value = 2356.4
$ws.Range("J112").Value = 8624.714
$ws.Range("K112").Value = 7069.200000000001
$ws.Range("L112").Value = 25874.142
$ws.Range("M112").Value = -5961.200000000001

$ws.Range("H129").Value = 753
$ws.Range("I129").Value = 753
$ws.Range("J129").Value = 0
$ws.Range("K129").Value = 2259
$ws.Range("L129").Value = 0
$ws.Range("M129").Value = 2741

$ws.Range("H138").Value = 2174.0657
$ws.Range("I138").Value = 1705.0938
$ws.Range("J138").Value = 2691.5518
$ws.Range("K138").Value = 5115.2814
$ws.Range("L138").Value = 8074.655400000001
$ws.Range("M138").Value = 24.71860000000015
$ws.Range("N138").Value = -18354.6554

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9920.566999999999
$ws.Range("I32").Value = 4873.4517
$ws.Range("J32").Value = 35997.332
$ws.Range("K32").Value = 4873.4517
$ws.Range("L32").Value = 35997.332
$ws.Range("M32").Value = -4586.4517

$ws.Range("H45").Value = 2771.182
$ws.Range("I45").Value = 2414.9443
$ws.Range("J45").Value = 4374.25
$ws.Range("K45").Value = 2414.9443
$ws.Range("L45").Value = 4374.25
$ws.Range("M45").Value = -2037.9443

$ws.Range("H132").Value = 1940.0857
$ws.Range("I132").Value = 1265.4073
$ws.Range("J132").Value = 4217.125
$ws.Range("K132").Value = 3796.2219
$ws.Range("L132").Value = 12651.375
$ws.Range("M132").Value = -1266.2219
$ws.Range("N132").Value = -17711.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2222
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 2222
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 2222
$ws.Range("N94").Value = -3124
$ws.Range("M94").ClearContents()

$ws.Range("H107").Value = 1906.5834
$ws.Range("I107").Value = 1238.2
$ws.Range("J107").Value = 2384
$ws.Range("K107").Value = 1238.2
$ws.Range("L107").Value = 2384
$ws.Range("M107").Value = 681.8
$ws.Range("N107").Value = -6224

$ws.Range("H132").Value = 142463.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 142463.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 142463.25
$ws.Range("N132").Value = -152583.25

$ws.Range("H134").Value = 4541.7856
$ws.Range("I134").Value = 2998.8057
$ws.Range("J134").Value = 13799.667
$ws.Range("K134").Value = 8996.417099999999
$ws.Range("L134").Value = 41399.001
$ws.Range("M134").Value = -6461.417099999999
$ws.Range("N134").Value = -46469.001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 7068
$ws.Range("I2").Value = 602
$ws.Range("J2").Value = 20000
$ws.Range("K2").Value = 602
$ws.Range("L2").Value = 20000
$ws.Range("M2").Value = -489

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 0
$ws.Range("M17").ClearContents()

$ws.Range("H31").Value = 7605.0244
$ws.Range("I31").Value = 3681.7273
$ws.Range("J31").Value = 12147.789
$ws.Range("K31").Value = 3681.7273
$ws.Range("L31").Value = 12147.789
$ws.Range("M31").Value = -3386.7273
$ws.Range("N31").Value = -12737.789

$ws.Range("H34").Value = 7605.0244
$ws.Range("I34").Value = 3681.7273
$ws.Range("J34").Value = 12147.789
$ws.Range("K34").Value = 3681.7273
$ws.Range("L34").Value = 12147.789
$ws.Range("M34").Value = -3479.7273
$ws.Range("N34").Value = -12551.789

$ws.Range("H132").Value = 2191.3225
$ws.Range("I132").Value = 1355.9434
$ws.Range("J132").Value = 7110.778
$ws.Range("K132").Value = 4067.8302
$ws.Range("L132").Value = 21332.334
$ws.Range("M132").Value = -1537.8302

$ws.Range("H141").Value = 202465.16
$ws.Range("I141").Value = 56799
$ws.Range("J141").Value = 241309.47
$ws.Range("K141").Value = 56799
$ws.Range("L141").Value = 241309.47
$ws.Range("M141").Value = -51619
$ws.Range("N141").Value = -251669.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1306.2
$ws.Range("I5").Value = 1194.3334
$ws.Range("J5").Value = 1474
$ws.Range("K5").Value = 3583.0002
$ws.Range("L5").Value = 4422
$ws.Range("M5").Value = -3471.0002

$ws.Range("H32").Value = 2847
$ws.Range("I32").Value = 2234
$ws.Range("J32").Value = 3000.25
$ws.Range("K32").Value = 6702
$ws.Range("L32").Value = 9000.75
$ws.Range("M32").Value = -6419
$ws.Range("N32").Value = -9566.75

$ws.Range("H46").Value = 4570.143
$ws.Range("I46").Value = 3999
$ws.Range("J46").Value = 4665.3335
$ws.Range("K46").Value = 11997
$ws.Range("L46").Value = 13996.0005
$ws.Range("M46").Value = -11906
$ws.Range("N46").Value = -14178.0005

$ws.Range("H121").Value = 904985.4399999999
$ws.Range("I121").Value = 192.5
$ws.Range("J121").Value = 1055784.2
$ws.Range("K121").Value = 577.5
$ws.Range("L121").Value = 3167352.6
$ws.Range("M121").Value = 732.5
$ws.Range("N121").Value = -3169972.6

$ws.Range("H135").Value = 1306.2
$ws.Range("I135").Value = 1194.3334
$ws.Range("J135").Value = 1474
$ws.Range("K135").Value = 10749.0006
$ws.Range("L135").Value = 13266
$ws.Range("M135").Value = -8214.000599999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7535
$ws.Range("I70").Value = 6919
$ws.Range("J70").Value = 9999
$ws.Range("K70").Value = 6919
$ws.Range("L70").Value = 9999
$ws.Range("M70").Value = -6649
$ws.Range("N70").Value = -10539

$ws.Range("H73").Value = 7535
$ws.Range("I73").Value = 6919
$ws.Range("J73").Value = 9999
$ws.Range("K73").Value = 6919
$ws.Range("L73").Value = 9999
$ws.Range("M73").Value = -5983
$ws.Range("N73").Value = -11871

$ws.Range("H80").Value = 3602.5
$ws.Range("I80").Value = 2205
$ws.Range("J80").Value = 5000
$ws.Range("K80").Value = 2205
$ws.Range("L80").Value = 5000
$ws.Range("M80").Value = -1207
$ws.Range("N80").Value = -6996

$ws.Range("H83").Value = 3602.5
$ws.Range("I83").Value = 2205
$ws.Range("J83").Value = 5000
$ws.Range("K83").Value = 11025
$ws.Range("L83").Value = 25000
$ws.Range("M83").Value = -6033
$ws.Range("N83").Value = -34984

$ws.Range("H102").Value = 1951
$ws.Range("I102").Value = 1132.4584
$ws.Range("J102").Value = 4133.778
$ws.Range("K102").Value = 1132.4584
$ws.Range("L102").Value = 4133.778
$ws.Range("M102").Value = 489.5416

$ws.Range("H122").Value = 4606.3447
$ws.Range("I122").Value = 2247.5334
$ws.Range("J122").Value = 7133.643
$ws.Range("K122").Value = 6742.600199999999
$ws.Range("L122").Value = 21400.929
$ws.Range("M122").Value = -4292.600199999999
$ws.Range("N122").Value = -26300.929

$ws.Range("H126").Value = 7415
$ws.Range("I126").Value = 7591.7
$ws.Range("J126").Value = 6973.25
$ws.Range("K126").Value = 22775.1
$ws.Range("L126").Value = 20919.75
$ws.Range("M126").Value = -20305.1

$ws.Range("H132").Value = 4807.6484
$ws.Range("I132").Value = 3096.842
$ws.Range("J132").Value = 6613.5
$ws.Range("K132").Value = 9290.526
$ws.Range("L132").Value = 19840.5
$ws.Range("M132").Value = -6760.526
$ws.Range("N132").Value = -24900.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 40000
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 40000
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 40000
$ws.Range("N42").Value = -41126
$ws.Range("M42").ClearContents()

$ws.Range("H49").Value = 40000
$ws.Range("I49").Value = 0
$ws.Range("J49").Value = 40000
$ws.Range("K49").Value = 0
$ws.Range("L49").Value = 40000
$ws.Range("N49").Value = -40294
$ws.Range("M49").ClearContents()

$ws.Range("H68").Value = 6041.875
$ws.Range("I68").Value = 3855.0625
$ws.Range("J68").Value = 10415.5
$ws.Range("K68").Value = 3855.0625
$ws.Range("L68").Value = 10415.5
$ws.Range("M68").Value = -3106.0625
$ws.Range("N68").Value = -11913.5

$ws.Range("H71").Value = 6041.875
$ws.Range("I71").Value = 3855.0625
$ws.Range("J71").Value = 10415.5
$ws.Range("K71").Value = 19275.3125
$ws.Range("L71").Value = 52077.5
$ws.Range("M71").Value = -15531.3125
$ws.Range("N71").Value = -59565.5

$ws.Range("H132").Value = 5663.5093
$ws.Range("I132").Value = 4379.8057
$ws.Range("J132").Value = 8381.941000000001
$ws.Range("K132").Value = 13139.4171
$ws.Range("L132").Value = 25145.823
$ws.Range("M132").Value = -10609.4171
$ws.Range("N132").Value = -30205.823

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1883.7192
$ws.Range("I132").Value = 1515.4222
$ws.Range("J132").Value = 3264.8333
$ws.Range("K132").Value = 4546.2666
$ws.Range("L132").Value = 9794.499899999999
$ws.Range("M132").Value = -2016.2666
$ws.Range("N132").Value = -14854.4999
